$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 1.55
$ws.Range("H3").Value = 4.5
$ws.Range("I3").Value = 5.25
$ws.Range("J3").Value = 2.05
$ws.Range("L3").Value = 5.5
$ws.Range("X3").Value = 8.5
$ws.Range("AJ3").Value = 17
$ws.Range("AL3").Value = 41
$ws.Range("AM3").Value = 41
$ws.Range("AO3").Value = 7.5
$ws.Range("AQ3").Value = 21
$ws.Range("AU3").Value = 8
$ws.Range("AV3").Value = 51
$ws.Range("BA3").Value = 101

# Row 4 updates
$ws.Range("Q4").Value = 2.15
$ws.Range("R4").Value = 1.67
